# MINNESOTA_2016.xlsx cleanup script
# - Rename header columns to machine-friendly names
# - Title-case Spanish connector words (de, del, el, la, los, las, y) within
#   municipality / state names throughout the data rows
# - Fix a floating point rounding artifact in D903
# - Remove the trailing footnote / metadata rows (1129-1134) and shrink the
#   used range accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header row: translate Spanish column headers to short machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case the Spanish connector words in columns A and B for all data rows.
#    Only whole-word, lowercase matches are capitalized (so names that already
#    legitimately start with a capitalized "El", e.g. "El Fuerte", stay untouched).
$capSet = @('de', 'del', 'el', 'la', 'los', 'las', 'y')
$lastRow = 1128

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -is [string] -and $v -ne "") {
            $words = $v.Split(' ')
            $changed = $false
            for ($i = 0; $i -lt $words.Length; $i++) {
                $w = $words[$i]
                if ($capSet -contains $w) {
                    $words[$i] = $w.Substring(0, 1).ToUpper() + $w.Substring(1)
                    $changed = $true
                }
            }
            if ($changed) {
                $cell.Value = [string]::Join(' ', $words)
            }
        }
    }
}

# 3. Fix the floating point rounding artifact on D903
$ws.Range("D903").Value = 0.009433962264150945

# 4. Remove the trailing footnote/metadata rows (1129 is blank, 1130-1134 hold
#    free text notes) so the sheet ends at row 1128.
$ws.Range("A1130:A1134").EntireRow.Delete()
